$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Ativação:" (row 8) date changes from 01/01/2018 to 01/01/2023.
# Assign via a literal-string formula (so Excel doesn't auto-convert the
# date-looking text into a real date serial number), then use
# Copy/PasteSpecial(values) to collapse the formula down to a plain cached
# text value so the cell keeps the same shared-string text type/style as
# before (no <f> left behind, no new number-format style created).
$xlPasteValues = -4163

$ws.Range("B8").Formula = '="01/01/2023"'
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial($xlPasteValues)
$ws.Range("C8").Formula = '="01/01/2023"'
$ws.Range("C8").Copy()
$ws.Range("C8").PasteSpecial($xlPasteValues)

# "Programa resumido:" (row 13) previously held "Semestral"; it now shows
# the same activation date value as row 8.
$ws.Range("B13").Formula = '="01/01/2023"'
$ws.Range("B13").Copy()
$ws.Range("B13").PasteSpecial($xlPasteValues)
$ws.Range("C13").Formula = '="01/01/2023"'
$ws.Range("C13").Copy()
$ws.Range("C13").PasteSpecial($xlPasteValues)

# "Programa:" (row 15) previously (incorrectly) echoed the activation date;
# it now shows the responsible professor's name.
$ws.Range("B15").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("C15").Value = "5840712 - Ângelo Capri Neto"

# "Método:" (row 18) professor changes from Ângelo Capri Neto to Rosa Ana Conte.
$ws.Range("B18").Value = "5840521 - Rosa Ana Conte"
$ws.Range("C18").Value = "5840521 - Rosa Ana Conte"
